# Commit: "slides - debug,vars, inventory vars, directory layouts"
# - Fix typo on row 16 ("Übungs" -> "Übung")
# - Add two new logged time entries (rows 17 and 18)
# - Add a total-cost formula in D2 (C2 * 20)
# - Leave the active selection on D6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 17 doesn't yet have number formats applied (style index 1 everywhere),
# so copy the formatting pattern from row 7 (Stunden/Datum/Zeitraum/Tasks
# styles: plain / date / plain / plain) before filling in the values.
$ws.Range("A7:D7").Copy($ws.Range("A17:D17"))

# Row 18 should mirror row 16's formatting (plain / date / time / plain).
$ws.Range("A16:D16").Copy($ws.Range("A18:D18"))

# New row 17: 13 Feb 2019
$ws.Range("A17").Value = 3.5
$d17 = Get-Date -Year 2019 -Month 2 -Day 13 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("B17").Value = $d17
$ws.Range("C17").Value = "20:00 - 23:30"
$ws.Range("D17").Value = "Präsentation/Übung 2 - simple playbook"

# Fix the existing typo on row 16 ("Übungs" -> "Übung")
$ws.Range("D16").Value = "Präsentation/Übung 1 - ansible_facts/aws fixes"

# New row 18: 14 Feb 2019
$ws.Range("A18").Value = 2
$d18 = Get-Date -Year 2019 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("B18").Value = $d18
$ws.Range("D18").Value = "Präsentation - Vars / Inventory Layouts / group_vars,host_vars"
$ws.Range("C18").Value = "21:00 - 23:00"

# Total cost column: hours * 20
$ws.Range("D2").Formula = "=C2*20"

$excel.Calculate()

# Matches the final active selection recorded in the workbook.
$ws.Range("D6").Select()
